$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 159: LeetCode 2785 - Sort Vowels in a String ---
$ws.Cells.Item(159, 1).Value = 2785
$ws.Cells.Item(159, 2).Value = "Sort Vowels in a String"
$ws.Cells.Item(159, 3).Value = "#string"
$ws.Cells.Item(159, 4).Value = "medium"
$ws.Cells.Item(159, 5).Value = 1
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 10
$ws.Cells.Item(159, 8).Value = 45912
$ws.Cells.Item(159, 8).NumberFormat = "m/d/yy"
$ws.Cells.Item(159, 9).Value = 45912
$ws.Cells.Item(159, 9).NumberFormat = "m/d/yy"
$ws.Rows.Item(159).RowHeight = 17

# --- Row 160: LeetCode 3227 - Vowels Game in a String ---
$ws.Cells.Item(160, 1).Value = 3227
$ws.Cells.Item(160, 2).Value = "Vowels Game in a String"
$ws.Cells.Item(160, 3).Value = "#math #string #greedy "
$ws.Cells.Item(160, 4).Value = "medium"
$ws.Cells.Item(160, 5).Value = 1
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 21
$ws.Cells.Item(160, 8).Value = 45913
$ws.Cells.Item(160, 8).NumberFormat = "m/d/yy"
$ws.Cells.Item(160, 9).Value = 45913
$ws.Cells.Item(160, 9).NumberFormat = "m/d/yy"
$ws.Rows.Item(160).RowHeight = 34

# --- Selection / window view bookkeeping ---
$ws.Range("G156").Select()

$excel.ActiveWindow.Width = 29100
$excel.ActiveWindow.Height = 14540
